$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (C3=0) is removed entirely.
$ws.Rows("3").Delete()

# Row 2: only Principle (C2=10000) survives; everything else in the row is cleared
# (values + the date/bool formatting that used to live on A2/G2).
$ws.Range("A2").Clear()
$ws.Range("B2").Clear()
$ws.Range("D2").Clear()
$ws.Range("E2").Clear()
$ws.Range("F2").Clear()
$ws.Range("G2").Clear()

# Column G (and its header cell G1) loses the date-format style that used
# to be applied to it.
$ws.Columns("G").ClearFormats()

# Column A's cached "best fit" width shrinks now that the long datetime
# value is gone (only the "Date" header remains) - approximate the
# resulting autofit width.
$ws.Columns("A").ColumnWidth = 4.29

# Selection returns to the top-left cell instead of the old G2 selection.
$null = $ws.Range("A1").Select()
